$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# 1) Typo fix: a stray comma used as a separator between person names was
#    meant to be a period ("," -> ".") in a handful of "Razon social" /
#    "Nombre Fantasia" shared strings.
# -----------------------------------------------------------------------
$ws.Range("E33").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F33").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E53").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F53").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"

$ws.Range("E52").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E54").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"

# -----------------------------------------------------------------------
# 2) "Importe" column (H): the scraper used to emit amounts formatted the
#    Spanish/Argentine way ("1.234,56" — "." thousands separator, ","
#    decimal separator). The fix re-emits them as plain floating point
#    text ("1234.56") instead.
#
#    These amounts are stored as literal TEXT in the sheet (not as real
#    numbers), so we force the cells to Text format first - otherwise
#    Excel "helpfully" re-parses a string like "500.00" as the number 500
#    and drops the formatting we are trying to write. ClearFormats()
#    afterwards restores the cells' original (default) formatting, since
#    only the textual content changed, not the look of the column.
# -----------------------------------------------------------------------
$ws.Range("H2:H98").NumberFormat = "@"

$ws.Range("H2").Value = "500.00"
$ws.Range("H3").Value = "236.00"
$ws.Range("H4").Value = "112497.50"
$ws.Range("H5").Value = "76.00"
$ws.Range("H6").Value = "11453.00"
$ws.Range("H7").Value = "70000.00"
$ws.Range("H8").Value = "36300.00"
$ws.Range("H9").Value = "24.38"
$ws.Range("H10").Value = "17786.46"
$ws.Range("H11").Value = "21599.98"
$ws.Range("H12").Value = "49183.20"
$ws.Range("H13").Value = "536.40"
$ws.Range("H14").Value = "22280.97"
$ws.Range("H15").Value = "4425.85"
$ws.Range("H16").Value = "469.20"
$ws.Range("H17").Value = "6843.77"
$ws.Range("H18").Value = "689.52"
$ws.Range("H19").Value = "313.32"
$ws.Range("H20").Value = "213.00"
$ws.Range("H21").Value = "33.74"
$ws.Range("H22").Value = "699.00"
$ws.Range("H23").Value = "222.07"
$ws.Range("H24").Value = "151.24"
$ws.Range("H25").Value = "2126.56"
$ws.Range("H26").Value = "474.36"
$ws.Range("H27").Value = "196.00"
$ws.Range("H28").Value = "1027.60"
$ws.Range("H29").Value = "69.00"
$ws.Range("H30").Value = "25200.00"
$ws.Range("H31").Value = "826.69"
$ws.Range("H32").Value = "10.50"
$ws.Range("H33").Value = "22.22"
$ws.Range("H34").Value = "57.68"
$ws.Range("H35").Value = "4214.05"
$ws.Range("H36").Value = "205.00"
$ws.Range("H37").Value = "2416.00"
$ws.Range("H38").Value = "3192.00"
$ws.Range("H39").Value = "2320.00"
$ws.Range("H40").Value = "157.50"
$ws.Range("H41").Value = "207.00"
$ws.Range("H42").Value = "6600.00"
$ws.Range("H43").Value = "180.00"
$ws.Range("H44").Value = "219.58"
$ws.Range("H45").Value = "1102.00"
$ws.Range("H46").Value = "140.00"
$ws.Range("H47").Value = "1300.00"
$ws.Range("H48").Value = "12117.40"
$ws.Range("H49").Value = "1170.00"
$ws.Range("H50").Value = "3972.00"
$ws.Range("H51").Value = "778.00"
$ws.Range("H52").Value = "290.00"
$ws.Range("H53").Value = "330.20"
$ws.Range("H54").Value = "1095.00"
$ws.Range("H55").Value = "18.00"
$ws.Range("H56").Value = "130.00"
$ws.Range("H57").Value = "680.00"
$ws.Range("H58").Value = "653.80"
$ws.Range("H59").Value = "370.00"
$ws.Range("H60").Value = "4693.00"
$ws.Range("H61").Value = "834.75"
$ws.Range("H62").Value = "20.00"
$ws.Range("H63").Value = "173.71"
$ws.Range("H64").Value = "330.00"
$ws.Range("H65").Value = "126.05"
$ws.Range("H66").Value = "590.00"
$ws.Range("H67").Value = "9695.90"
$ws.Range("H68").Value = "204.40"
$ws.Range("H69").Value = "165.00"
$ws.Range("H70").Value = "552.50"
$ws.Range("H71").Value = "200.00"
$ws.Range("H72").Value = "700.00"
$ws.Range("H73").Value = "300.00"
$ws.Range("H74").Value = "1815.00"
$ws.Range("H75").Value = "435.00"
$ws.Range("H76").Value = "1424.83"
$ws.Range("H77").Value = "600.00"
$ws.Range("H78").Value = "300.00"
$ws.Range("H79").Value = "750.00"
$ws.Range("H80").Value = "120.00"
$ws.Range("H81").Value = "8520.50"
$ws.Range("H82").Value = "339.30"
$ws.Range("H83").Value = "336.00"
$ws.Range("H84").Value = "6656.95"
$ws.Range("H85").Value = "439.04"
$ws.Range("H86").Value = "19923.00"
$ws.Range("H87").Value = "477.00"
$ws.Range("H88").Value = "917.42"
$ws.Range("H89").Value = "7539.86"
$ws.Range("H90").Value = "2842.00"
$ws.Range("H91").Value = "4192.50"
$ws.Range("H92").Value = "896.59"
$ws.Range("H93").Value = "2965.71"
$ws.Range("H94").Value = "944.16"
$ws.Range("H95").Value = "403840.04"
$ws.Range("H96").Value = "81000.00"
$ws.Range("H97").Value = "56.00"
$ws.Range("H98").Value = "894.00"

$ws.Range("H2:H98").ClearFormats()
